# My_Trading_Bot_Journey.pptx — "Updated PPTX Journey Deck"
#
# Re-themes the deck from "Building My First Python Trading Bot" to
# "The Developer's Journey" (Building & Publishing a Python Trading Bot /
# From Zero to GitHub), rewrites the four existing step slides, and
# appends three new slides (Version Control, Publishing to GitHub, and a
# replacement "Future Roadmap" slide).

$p = $ppt.ActivePresentation

function Set-Body($shape, [string[]]$paragraphs) {
    # Setting TextRange.Text directly onto an existing multi-paragraph
    # run tends to diff the old/new strings paragraph-by-paragraph and
    # split runs on any surviving common prefix/suffix. Collapsing the
    # frame down to a single throw-away paragraph first guarantees the
    # following assignment lands as one clean run per paragraph.
    $tr = $shape.TextFrame.TextRange
    $tr.Text = "~"
    $tr.Text = [string]::Join("`r", $paragraphs)
}

# ---------------------------------------------------------------------
# Slide 1 — Title slide
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
Set-Body $s1.Shapes.Item(1) @("The Developer's Journey")
Set-Body $s1.Shapes.Item(2) @(
    "Building & Publishing a Python Trading Bot",
    "From Zero to GitHub"
)

# ---------------------------------------------------------------------
# Slide 2 — Phase 1: The Setup -> Step 1: The Vision
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
Set-Body $s2.Shapes.Item(1) @("Step 1: The Vision")
Set-Body $s2.Shapes.Item(2) @(
    "Objective: Build an automated stock trading assistant.",
    "The Problem: We cannot watch the screen 24/7.",
    "The Solution: A Python script that watches, thinks, and alerts.",
    "Tools Selected: VS Code (Editor), Python 3.13 (Engine)."
)

# ---------------------------------------------------------------------
# Slide 3 — Phase 2: The First Script -> Step 2: The Logic (The 'Trifecta')
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
Set-Body $s3.Shapes.Item(1) @("Step 2: The Logic (The 'Trifecta')")
Set-Body $s3.Shapes.Item(2) @(
    "We programmed the bot to think like a disciplined trader.",
    "Indicator 1: SMA (20) - Is the trend up?",
    "Indicator 2: RSI (14) - Is the price fair?",
    "Indicator 3: MACD (12,26,9) - Is momentum building?",
    "Rule: The bot only signals when ALL three agree."
)

# ---------------------------------------------------------------------
# Slide 4 — Phase 3: The Loop (Bot Mode) -> Step 3: Automation & Alerts
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
Set-Body $s4.Shapes.Item(1) @("Step 3: Automation & Alerts")
Set-Body $s4.Shapes.Item(2) @(
    "Challenge: Making it run forever.",
    "Solution: Implemented a 'While True' infinite loop.",
    "The Voice: Added 'plyer' to trigger Windows Desktop Notifications.",
    "The Result: The user can work on other tasks while the bot guards the portfolio."
)

# ---------------------------------------------------------------------
# Slide 5 — Phase 4: Future Goals -> Step 4: Persistent Logging
# ---------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
Set-Body $s5.Shapes.Item(1) @("Step 4: Persistent Logging")
Set-Body $s5.Shapes.Item(2) @(
    "Challenge: How do we know what the bot did while we slept?",
    "Solution: Added a CSV Logger.",
    "Mechanism: Python opens 'trade_log.csv' and appends every Buy/Sell decision.",
    "Benefit: Creates a permanent audit trail for backtesting and review."
)

# ---------------------------------------------------------------------
# Slide 6 (new) — Step 5: Version Control (Git)
# ---------------------------------------------------------------------
$s6 = $p.Slides.Add(6, 2)
Set-Body $s6.Shapes.Item(1) @("Step 5: Version Control (Git)")
Set-Body $s6.Shapes.Item(2) @(
    "We turned our folder into a Repository.",
    "Git Init: Started tracking changes.",
    "Git Add/Commit: Saved 'snapshots' of our code.",
    "Why? To prevent data loss and allow us to 'rewind' if we break something."
)

# ---------------------------------------------------------------------
# Slide 7 (new) — Step 6: Publishing to the Cloud (GitHub)
# ---------------------------------------------------------------------
$s7 = $p.Slides.Add(7, 2)
Set-Body $s7.Shapes.Item(1) @("Step 6: Publishing to the Cloud (GitHub)")
Set-Body $s7.Shapes.Item(2) @(
    "Final Step: Uploading to the world.",
    "Remote Repo: Created a secure box on GitHub.com.",
    "Push: Sent our local code to the cloud.",
    "Authentication: Secured the connection using a Personal Access Token (PAT).",
    "Outcome: The code is now safe, shareable, and professional."
)

# ---------------------------------------------------------------------
# Slide 8 (new) — Future Roadmap
# ---------------------------------------------------------------------
$s8 = $p.Slides.Add(8, 2)
Set-Body $s8.Shapes.Item(1) @("Future Roadmap")
Set-Body $s8.Shapes.Item(2) @(
    "Phase 2: Cloud Hosting (Running 24/7 on a server).",
    "Phase 3: Salesforce Integration (Logging trades as CRM records).",
    "Phase 4: Backtesting (Simulating performance on past data).",
    "Status: Phase 1 Complete."
)
